$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "supplydemand" to "SupplyDemand"
$ws.Name = "SupplyDemand"

# Normalize the selection to just B4 (was previously a full row/column select)
$ws.Range("B4").Select() | Out-Null
